$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.02729999999999
$ws.Range("C4").Value = -12.5774
$ws.Range("B7").Value = 4.966299999999996
$ws.Range("A8").Value = -22.34270000000001
$ws.Range("A10").Value = -21.69239999999999
$ws.Range("E10").Value = 16.47389999999999
$ws.Range("C11").Value = -11.5915
$ws.Range("A12").Value = -21.61280000000001
$ws.Range("E12").Value = 17.8538
$ws.Range("E13").Value = 16.3765
$ws.Range("B14").Value = 4.809599999999993
$ws.Range("C14").Value = -13.5295
$ws.Range("E14").Value = 16.43210000000001
$ws.Range("B15").Value = 4.608699999999994
$ws.Range("A18").Value = -21.9537
$ws.Range("B18").Value = 5.943800000000001
$ws.Range("C18").Value = -11.91710000000001
$ws.Range("C19").Value = -11.17760000000001
$ws.Range("B20").Value = 8.575700000000001
$ws.Range("C21").Value = -12.0145
$ws.Range("A25").Value = -21.49629999999999
$ws.Range("C27").Value = -12.9807
$ws.Range("B29").Value = 5.398400000000004
$ws.Range("E29").Value = 17.36090000000001
$ws.Range("B30").Value = 4.813600000000002
$ws.Range("B31").Value = 5.248000000000002
$ws.Range("C31").Value = -12.9982
$ws.Range("E32").Value = 16.47359999999999
$ws.Range("B35").Value = 8.439000000000004
$ws.Range("E35").Value = 16.31840000000001
$ws.Range("A37").Value = -18.97029999999999
$ws.Range("C38").Value = -12.8747
$ws.Range("B40").Value = 9.12099999999999
$ws.Range("C42").Value = -12.1231
$ws.Range("E43").Value = 17.24950000000002
$ws.Range("B44").Value = 4.7307
$ws.Range("C44").Value = -13.6455
$ws.Range("C47").Value = -12.6351
$ws.Range("E48").Value = 17.42180000000003
$ws.Range("E49").Value = 15.87359999999999
$ws.Range("B50").Value = 4.2866
$ws.Range("E50").Value = 16.27609999999999
$ws.Range("E51").Value = 17.09880000000001
$ws.Range("B54").Value = 4.589299999999998
$ws.Range("A55").Value = -21.8353
$ws.Range("C56").Value = -13.69829999999999
$ws.Range("E56").Value = 16.5261
$ws.Range("C58").Value = -12.5103
$ws.Range("E61").Value = 16.71210000000001
$ws.Range("C65").Value = -12.1773
$ws.Range("A68").Value = -21.46510000000001
$ws.Range("B68").Value = 4.533499999999997
$ws.Range("E69").Value = 17.34240000000004
$ws.Range("E71").Value = 16.90580000000001
$ws.Range("C73").Value = -13.08110000000001
$ws.Range("B76").Value = 6.037499999999997
$ws.Range("A77").Value = -20.22329999999998
$ws.Range("A78").Value = -20.17299999999997
$ws.Range("A79").Value = -20.03119999999998
$ws.Range("E79").Value = 18.07650000000002
$ws.Range("A80").Value = -19.68919999999997
$ws.Range("A81").Value = -21.78420000000001
$ws.Range("E81").Value = 16.74389999999999
$ws.Range("A82").Value = -21.7927
$ws.Range("A84").Value = -21.92620000000001
$ws.Range("B87").Value = 4.670899999999996
$ws.Range("B88").Value = 4.738299999999998
$ws.Range("C90").Value = -13.4143
$ws.Range("B92").Value = 5.343399999999994
$ws.Range("C92").Value = -10.51509999999999
$ws.Range("E92").Value = 18.24570000000003
$ws.Range("C94").Value = -10.0064
$ws.Range("C95").Value = -12.3487
$ws.Range("B96").Value = 4.692900000000006
$ws.Range("B98").Value = 6.773399999999999
$ws.Range("A101").Value = -21.71579999999999
$ws.Range("B101").Value = 5.5914
$ws.Range("C101").Value = -12.6284
$ws.Range("A102").Value = -21.75380000000001
$ws.Range("B102").Value = 5.056100000000003
